# 更新PPT 06 07 08
# Applies the four text edits described by the commit:
#   1. Slide 13 title: merge "以及" + "原型链的图解" -> "以及原型链的图解"
#   2. Slide  5 body : "语言继承方式简介" -> "语言继承方式"
#   3. Slide  5 title: "JS" / "对象简介" -> "JS" / "对象及继承方式综述"
#   4. Slide  8 note : merge "后" + "半部分 属性相关操作" -> "后半部分 属性相关操作"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 13 - "基于构造函数实现的原型继承" + "以及" + "原型链的图解"
#    The last two runs get merged into a single run with the combined text,
#    leaving the first run (and all formatting) untouched.
# ---------------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)
$titleShape13 = $slide13.Shapes.Item(2)
$tr13 = $titleShape13.TextFrame.TextRange
$tr13.Characters(14, 8).Text = "以及原型链的图解"

# ---------------------------------------------------------------------------
# 2) Slide 5 - content placeholder: "JavaScript" + "语言继承方式简介" -> 去掉“简介”
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$bodyShape5 = $slide5.Shapes.Item(1)
$trBody5 = $bodyShape5.TextFrame.TextRange
$trBody5.Characters(11, 8).Text = "语言继承方式"

# ---------------------------------------------------------------------------
# 3) Slide 5 - title placeholder: "JS" + "对象简介" -> "JS" + "对象及继承方式综述"
#    Keep "JS" as its own run and only rewrite the Chinese run's text.
# ---------------------------------------------------------------------------
$titleShape5 = $slide5.Shapes.Item(2)
$trTitle5 = $titleShape5.TextFrame.TextRange
$trTitle5.Characters(1, 2).Text = "JS"
$trTitle5.Characters(3, 4).Text = "对象及继承方式综述"

# ---------------------------------------------------------------------------
# 4) Slide 8 - last textbox: "参见实例" + "demo03" + "后" + "半部分 属性相关操作"
#    Merge the last two runs into one.
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$noteShape8 = $slide8.Shapes.Item($slide8.Shapes.Count)
$tr8 = $noteShape8.TextFrame.TextRange
$tr8.Characters(11, 11).Text = "后半部分 属性相关操作"
